$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 284.55554
$ws.Range("I4").Value = 203.83333
$ws.Range("K4").Value = 203.83333
$ws.Range("M4").Value = -89.83332999999999
$ws.Range("H6").Value = 426.11765
$ws.Range("I6").Value = 369.1
$ws.Range("J6").Value = 507.57144
$ws.Range("K6").Value = 1107.3
$ws.Range("L6").Value = 1522.71432
$ws.Range("M6").Value = -995.3000000000002
$ws.Range("N6").Value = -1746.71432
$ws.Range("H19").Value = 3334.5264
$ws.Range("I19").Value = 3517.6667
$ws.Range("J19").Value = 3169.7
$ws.Range("K19").Value = 3517.6667
$ws.Range("L19").Value = 3169.7
$ws.Range("M19").Value = -3342.6667
$ws.Range("N19").Value = -3519.7
$ws.Range("H31").Value = 80.5
$ws.Range("I31").Value = 80.5
$ws.Range("K31").Value = 241.5
$ws.Range("M31").Value = -11.5
$ws.Range("H32").Value = 1886.125
$ws.Range("I32").Value = 1531.3334
$ws.Range("J32").Value = 2099
$ws.Range("K32").Value = 1531.3334
$ws.Range("L32").Value = 2099
$ws.Range("M32").Value = -1205.3334
$ws.Range("N32").Value = -2751
$ws.Range("H41").Value = 512.25
$ws.Range("I41").Value = 414
$ws.Range("K41").Value = 414
$ws.Range("M41").Value = 26
$ws.Range("H86").Value = 8585.571
$ws.Range("J86").Value = 9333
$ws.Range("L86").Value = 9333
$ws.Range("N86").Value = -11579
$ws.Range("H89").Value = 8585.571
$ws.Range("J89").Value = 9333
$ws.Range("L89").Value = 46665
$ws.Range("N89").Value = -57897
$ws.Range("H132").Value = 2609.28
$ws.Range("I132").Value = 2563.2942
$ws.Range("K132").Value = 7689.882599999999
$ws.Range("M132").Value = -5159.882599999999
$ws.Range("H137").Value = 1957.6666
$ws.Range("I137").Value = 1527.7142
$ws.Range("K137").Value = 4583.142599999999
$ws.Range("M137").Value = -2033.142599999999
$ws.Range("H141").Value = 2078.6
$ws.Range("I141").Value = 2078.6
$ws.Range("K141").Value = 6235.799999999999
$ws.Range("M141").Value = -1055.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 30050
$ws.Range("J44").Value = 30050
$ws.Range("L44").Value = 30050
$ws.Range("N44").Value = -31026
$ws.Range("H51").Value = 41000
$ws.Range("J51").Value = 41000
$ws.Range("L51").Value = 41000
$ws.Range("N51").Value = -42512
$ws.Range("H63").Value = 4306.7
$ws.Range("I63").Value = 2261.1667
$ws.Range("J63").Value = 7375
$ws.Range("K63").Value = 2261.1667
$ws.Range("L63").Value = 7375
$ws.Range("M63").Value = -1575.1667
$ws.Range("N63").Value = -8747
$ws.Range("H66").Value = 4306.7
$ws.Range("I66").Value = 2261.1667
$ws.Range("J66").Value = 7375
$ws.Range("K66").Value = 11305.8335
$ws.Range("L66").Value = 36875
$ws.Range("M66").Value = -7873.833500000001
$ws.Range("N66").Value = -43739
$ws.Range("H97").Value = 1012.3158
$ws.Range("I97").Value = 639.625
$ws.Range("K97").Value = 639.625
$ws.Range("M97").Value = -143.625
$ws.Range("H102").Value = 600
$ws.Range("I102").Value = 600
$ws.Range("K102").Value = 600
$ws.Range("M102").Value = 1022
$ws.Range("H107").Value = 14999
$ws.Range("J107").Value = 14999
$ws.Range("L107").Value = 14999
$ws.Range("N107").Value = -22679
$ws.Range("H109").Value = 47499.5
$ws.Range("J109").Value = 47499.5
$ws.Range("L109").Value = 47499.5
$ws.Range("N109").Value = -50273.5
$ws.Range("H132").Value = 1008.25
$ws.Range("I132").Value = 1034.4286
$ws.Range("J132").Value = 825
$ws.Range("K132").Value = 3103.2858
$ws.Range("L132").Value = 2475
$ws.Range("M132").Value = -573.2857999999997
$ws.Range("N132").Value = -7535

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 366.66666
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = -127
$ws.Range("H86").Value = 11125
$ws.Range("I86").Value = 9250
$ws.Range("J86").Value = 13000
$ws.Range("K86").Value = 9250
$ws.Range("L86").Value = 13000
$ws.Range("M86").Value = -8127
$ws.Range("N86").Value = -15246
$ws.Range("H89").Value = 11125
$ws.Range("I89").Value = 9250
$ws.Range("J89").Value = 13000
$ws.Range("K89").Value = 46250
$ws.Range("L89").Value = 65000
$ws.Range("M89").Value = -40634
$ws.Range("N89").Value = -76232
$ws.Range("H105").Value = 4184.4614
$ws.Range("I105").Value = 3889.9
$ws.Range("J105").Value = 5166.3335
$ws.Range("K105").Value = 3889.9
$ws.Range("L105").Value = 5166.3335
$ws.Range("M105").Value = -2142.9
$ws.Range("N105").Value = -8660.333500000001
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774
$ws.Range("H112").Value = 36156
$ws.Range("J112").Value = 36156
$ws.Range("L112").Value = 36156
$ws.Range("N112").Value = -39110
$ws.Range("H134").Value = 1693
$ws.Range("I134").Value = 1562.3
$ws.Range("K134").Value = 4686.9
$ws.Range("M134").Value = -2151.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 19199.9
$ws.Range("J41").Value = 20777.666
$ws.Range("L41").Value = 20777.666
$ws.Range("N41").Value = -21633.666
$ws.Range("H105").Value = 1688
$ws.Range("I105").Value = 1766.1666
$ws.Range("J105").Value = 1531.6666
$ws.Range("K105").Value = 1766.1666
$ws.Range("L105").Value = 1531.6666
$ws.Range("M105").Value = -19.16660000000002
$ws.Range("N105").Value = -5025.6666
$ws.Range("H141").Value = 361744
$ws.Range("J141").Value = 361744
$ws.Range("L141").Value = 361744
$ws.Range("N141").Value = -372104

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 375
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H80").Value = 7645.3335
$ws.Range("J80").Value = 10293.25
$ws.Range("L80").Value = 30879.75
$ws.Range("N80").Value = -32751.75
$ws.Range("H83").Value = 7645.3335
$ws.Range("J83").Value = 10293.25
$ws.Range("L83").Value = 92639.25
$ws.Range("N83").Value = -101999.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 319
$ws.Range("I16").Value = 319
$ws.Range("K16").Value = 319
$ws.Range("M16").Value = -149
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H68").Value = 2099.1428
$ws.Range("I68").Value = 2138.8
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 2138.8
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -1389.8
$ws.Range("N68").Value = -3498
$ws.Range("H71").Value = 2099.1428
$ws.Range("I71").Value = 2138.8
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 10694
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -6950
$ws.Range("N71").Value = -17488
$ws.Range("H93").Value = 1600
$ws.Range("I93").Value = 1450
$ws.Range("K93").Value = 1450
$ws.Range("M93").Value = -202
$ws.Range("H139").Value = 90000
$ws.Range("J139").Value = 90000
$ws.Range("L139").Value = 90000
$ws.Range("N139").Value = -100280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 44188
$ws.Range("J74").Value = 34497.5
$ws.Range("L74").Value = 34497.5
$ws.Range("N74").Value = -36369.5
$ws.Range("H77").Value = 44188
$ws.Range("J77").Value = 34497.5
$ws.Range("L77").Value = 103492.5
$ws.Range("N77").Value = -112852.5
$ws.Range("H122").Value = 2292.5715
$ws.Range("I122").Value = 2009.1
$ws.Range("J122").Value = 3001.25
$ws.Range("K122").Value = 6027.299999999999
$ws.Range("L122").Value = 9003.75
$ws.Range("M122").Value = -3577.299999999999
$ws.Range("N122").Value = -13903.75
$ws.Range("H139").Value = 80000
$ws.Range("J139").Value = 80000
$ws.Range("L139").Value = 80000
$ws.Range("N139").Value = -90280
